# Captura de pantalla y carga en oneDrive
# Update the "Constants" sheet: rename the ExExcelDirectoryPath / data_excel
# config entry to sucessFormDirectoryPath (Key and Value now match).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Range("A5").Value = "sucessFormDirectoryPath"
$ws.Range("B5").Value = "sucessFormDirectoryPath"
